# Added 4wk low sales check: recalculated forecast/coverage/seasonality
# numbers on "Forecast Comparison" and rolled the new totals up into
# "Summary".

$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet -------------------------------------------

# Row 2 (W10)
$wsForecast.Range("L2").Value = 1.17

# Row 3 (W11)
$wsForecast.Range("D3").Value = 19
$wsForecast.Range("H3").Value = 4.53
$wsForecast.Range("L3").Value = 0.89

# Row 4 (W12)
$wsForecast.Range("D4").Value = 19
$wsForecast.Range("H4").Value = 3.53
$wsForecast.Range("L4").Value = 1.18

# Row 5 (W13)
$wsForecast.Range("D5").Value = 19
$wsForecast.Range("H5").Value = 2.53
$wsForecast.Range("L5").Value = 1.17

# Row 6 (W14)
$wsForecast.Range("D6").Value = 19
$wsForecast.Range("H6").Value = 1.53
$wsForecast.Range("L6").Value = 1.1

# Row 7 (W15)
$wsForecast.Range("D7").Value = 19
$wsForecast.Range("H7").Value = 0.53
$wsForecast.Range("J7").Value = "Urgent"
$wsForecast.Range("L7").Value = 0.9399999999999999

# Row 8 (W16)
$wsForecast.Range("D8").Value = 20
$wsForecast.Range("H8").Value = 0
$wsForecast.Range("L8").Value = 1.01

# Row 9 (W17)
$wsForecast.Range("D9").Value = 21
$wsForecast.Range("L9").Value = 1.14

# Row 10 (W18)
$wsForecast.Range("D10").Value = 22
$wsForecast.Range("L10").Value = 1.03

# Row 11 (W19)
$wsForecast.Range("D11").Value = 23
$wsForecast.Range("L11").Value = 1.15

# Row 12 (W20)
$wsForecast.Range("D12").Value = 23
$wsForecast.Range("L12").Value = 0.84

# Row 13 (W21)
$wsForecast.Range("D13").Value = 22
$wsForecast.Range("L13").Value = 1.07

# Row 14 (W22)
$wsForecast.Range("D14").Value = 22
$wsForecast.Range("L14").Value = 1.14

# Row 15 (W23)
$wsForecast.Range("D15").Value = 23
$wsForecast.Range("L15").Value = 1.03

# Row 16 (W24)
$wsForecast.Range("D16").Value = 24
$wsForecast.Range("L16").Value = 1.1

# Row 17 (W25)
$wsForecast.Range("D17").Value = 24
$wsForecast.Range("L17").Value = 1.15

# --- Summary sheet ---------------------------------------------------------

$wsSummary.Range("B9").Value = "338"
$wsSummary.Range("B10").Value = "155"
$wsSummary.Range("B11").Value = "76"
$wsSummary.Range("B12").Value = "24"
$wsSummary.Range("B14").Value = "19"
